$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 - Sahara Telecom (was Mamun Telecom)
$ws.Range("B2").Value = "DSR-0350"
$ws.Range("C2").Value = "Sahara Telecom"
$ws.Range("D2").Value = "Lalpur"
$ws.Range("E2").Value = "Md Shahin Ali"
$ws.Range("I2").Value = "Md Shahin Ali"
$ws.Range("J2").Value = 1712617115
$ws.Range("K2").Value = "Natore"
$ws.Range("L2").Value = "Lalpur"
$ws.Range("N2").Value = "Abdulpur, Lalpur, Natore."
$ws.Range("P2").Value = 1712617115
$ws.Range("T2").Value = 1712617115

# Row 3 - Square  Telecom (was Masum Electronics)
$ws.Range("B3").Value = "DSR-0350"
$ws.Range("C3").Value = "Square  Telecom"
$ws.Range("D3").Value = "Lalpur"
$ws.Range("E3").Value = "Md Rony Ali"
$ws.Range("I3").Value = "Md Rony Ali"
$ws.Range("J3").Value = 1714588737
$ws.Range("K3").Value = "Natore"
$ws.Range("L3").Value = "Lalpur"
$ws.Range("N3").Value = "Gopalpur, Lalpur, "
$ws.Range("P3").Value = 1714588737
$ws.Range("T3").Value = 1714588737

# Row 4 - Sheikh Telecom 2 (was Moom Telecom)
$ws.Range("B4").Value = "DSR-0619"
$ws.Range("C4").Value = "Sheikh Telecom 2"
$ws.Range("D4").Value = "Bonpara"
$ws.Range("E4").Value = "Sheikh Saifuddin"
$ws.Range("I4").Value = "Sheikh Saifuddin"
$ws.Range("J4").Value = 1712337781
$ws.Range("K4").Value = "Natore"
$ws.Range("L4").Value = "Baraigram"
$ws.Range("M4").Value = "ZSO-0022"
$ws.Range("N4").Value = "Bonpara, Natore."
$ws.Range("P4").Value = 1712337781
$ws.Range("T4").Value = 1712337781

# Row 5 - Azim Mobile Center (fills the pre-existing blank row 5)
$ws.Range("A5").Value = "DEL-0179"
$ws.Range("B5").Value = "DSR-0619"
$ws.Range("C5").Value = "Azim Mobile Center"
$ws.Range("D5").Value = "Bonpara"
$ws.Range("E5").Value = "Md Azim Uddin"
$ws.Range("G5").Value = "GO"
$ws.Range("I5").Value = "Md Azim Uddin"
$ws.Range("J5").Value = 1760264390
$ws.Range("K5").Value = "Natore"
$ws.Range("L5").Value = "Baraigram"
$ws.Range("M5").Value = "ZSO-0022"
$ws.Range("N5").Value = "Koenbazar, Baraigram, Natore."
$ws.Range("P5").Value = 1760264390
$ws.Range("Q5").Value = "C"
$ws.Range("R5").Value = "Rural"
$ws.Range("S5").Value = "bKash"
$ws.Range("T5").Value = 1760264390

# Row 6 - RS Mobile (fills the pre-existing blank row 6)
$ws.Range("A6").Value = "DEL-0179"
$ws.Range("B6").Value = "DSR-0619"
$ws.Range("C6").Value = "RS Mobile "
$ws.Range("D6").Value = "Bonpara"
$ws.Range("E6").Value = "Md Rakib Ali"
$ws.Range("G6").Value = "GO"
$ws.Range("I6").Value = "Md Rakib Ali"
$ws.Range("J6").Value = 1710372747
$ws.Range("K6").Value = "Natore"
$ws.Range("L6").Value = "Baraigram"
$ws.Range("M6").Value = "ZSO-0022"
$ws.Range("N6").Value = "Koenbazar, Baraigram, Natore."
$ws.Range("P6").Value = 1710372747
$ws.Range("Q6").Value = "C"
$ws.Range("R6").Value = "Rural"
$ws.Range("S6").Value = "bKash"
$ws.Range("T6").Value = 1710372747

# Row 7 - Gourango Hardware (fills the pre-existing blank row 7)
$ws.Range("A7").Value = "DEL-0179"
$ws.Range("B7").Value = "DSR-0619"
$ws.Range("C7").Value = "Gourango Hardware"
$ws.Range("D7").Value = "Bonpara"
$ws.Range("E7").Value = "Sree Gones Chandro"
$ws.Range("G7").Value = "GO"
$ws.Range("I7").Value = "Sree Gones Chandro"
$ws.Range("J7").Value = 1761689867
$ws.Range("K7").Value = "Natore"
$ws.Range("L7").Value = "Baraigram"
$ws.Range("M7").Value = "ZSO-0022"
$ws.Range("N7").Value = "Shahebbazar, Baraigram, Natore."
$ws.Range("P7").Value = 1761689867
$ws.Range("Q7").Value = "C"
$ws.Range("R7").Value = "Rural"
$ws.Range("S7").Value = "bKash"
$ws.Range("T7").Value = 1761689867

# Row 8 - Nirob Mobile & Computer (fills the pre-existing blank row 8)
$ws.Range("A8").Value = "DEL-0179"
$ws.Range("B8").Value = "DSR-0619"
$ws.Range("C8").Value = "Nirob Mobile & Computer"
$ws.Range("D8").Value = "Bonpara"
$ws.Range("E8").Value = "Md Nirob Mahabur Rahman"
$ws.Range("G8").Value = "GO"
$ws.Range("I8").Value = "Md Nirob Mahabur Rahman"
$ws.Range("J8").Value = 1791953259
$ws.Range("K8").Value = "Natore"
$ws.Range("L8").Value = "Bagatipara"
$ws.Range("M8").Value = "ZSO-0022"
$ws.Range("N8").Value = "Dayarampur, Bagatipara, Natore."
$ws.Range("P8").Value = 1791953259
$ws.Range("Q8").Value = "C"
$ws.Range("R8").Value = "Rural"
$ws.Range("S8").Value = "bKash"
$ws.Range("T8").Value = 1791953259

# Row 9 (previously all-blank row) picks up the bordered-cell styling pattern
# that row 8 used to have, to make room for the new row 10 below it.
$ws.Range("D9").Borders.Item(7).LineStyle = 1
$ws.Range("D9").Borders.Item(8).LineStyle = 1
$ws.Range("D9").Borders.Item(9).LineStyle = 1
$ws.Range("D9").Borders.Item(10).LineStyle = 1
$ws.Range("K9").Borders.Item(7).LineStyle = 1
$ws.Range("K9").Borders.Item(8).LineStyle = 1
$ws.Range("K9").Borders.Item(9).LineStyle = 1
$ws.Range("K9").Borders.Item(10).LineStyle = 1
$ws.Range("L9").Borders.Item(7).LineStyle = 1
$ws.Range("L9").Borders.Item(8).LineStyle = 1
$ws.Range("L9").Borders.Item(9).LineStyle = 1
$ws.Range("L9").Borders.Item(10).LineStyle = 1
$ws.Range("N9").Borders.Item(7).LineStyle = 1
$ws.Range("N9").Borders.Item(8).LineStyle = 1
$ws.Range("N9").Borders.Item(9).LineStyle = 1
$ws.Range("N9").Borders.Item(10).LineStyle = 1

# Fresh blank row 10 (matches the borderless style the old row 9 used to have)
$ws.Range("A10:T10").Value = ""

# The trailing "blank" marker row moves from row 20 to row 21
$ws.Range("E20").Value = ""
$ws.Range("E21").Value = "   "

# Column E widened slightly to fit the longer data
$ws.Columns("E:E").ColumnWidth = 25.7109375

# Update sheet dimension / selection to match the post-edit state
$ws.Range("G22").Select()
